$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17
$ws.Range("C2").Value = "dog/dog072.png"
$ws.Range("D2").Value = "gründen"
$ws.Range("E2").Value = "dog"
$ws.Range("B3").Value = 111
$ws.Range("C3").Value = "dog/dog121.png"
$ws.Range("D3").Value = "saufen"
$ws.Range("E3").Value = "dog"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = "dog/dog068.png"
$ws.Range("D4").Value = "sondern"
$ws.Range("E4").Value = "dog"
$ws.Range("B5").Value = 13
$ws.Range("C5").Value = "dog/dog069.png"
$ws.Range("D5").Value = "drehen"
$ws.Range("E5").Value = "dog"
$ws.Range("B6").Value = 100
$ws.Range("C6").Value = "flower/flower089.png"
$ws.Range("D6").Value = "segeln"
$ws.Range("E6").Value = "flower"
$ws.Range("B7").Value = 54
$ws.Range("C7").Value = "flower/flower101.png"
$ws.Range("D7").Value = "bleiben"
$ws.Range("E7").Value = "flower"
$ws.Range("B8").Value = 33
$ws.Range("C8").Value = "dog/dog083.png"
$ws.Range("D8").Value = "schmecken"
$ws.Range("E8").Value = "dog"
$ws.Range("B9").Value = 86
$ws.Range("C9").Value = "dog/dog070.png"
$ws.Range("D9").Value = "drohen"
$ws.Range("E9").Value = "dog"
$ws.Range("B10").Value = 52
$ws.Range("C10").Value = "flower/flower072.png"
$ws.Range("D10").Value = "lehnen"
$ws.Range("E10").Value = "flower"
$ws.Range("B11").Value = 21
$ws.Range("C11").Value = "flower/flower087.png"
$ws.Range("D11").Value = "stärken"
$ws.Range("E11").Value = "flower"
$ws.Range("B12").Value = 26
$ws.Range("C12").Value = "flower/flower091.png"
$ws.Range("D12").Value = "kehren"
$ws.Range("E12").Value = "flower"
$ws.Range("B13").Value = 72
$ws.Range("C13").Value = "flower/flower081.png"
$ws.Range("D13").Value = "fliegen"
$ws.Range("E13").Value = "flower"
$ws.Range("B14").Value = 98
$ws.Range("C14").Value = "dog/dog086.png"
$ws.Range("D14").Value = "haken"
$ws.Range("E14").Value = "dog"
$ws.Range("B15").Value = 37
$ws.Range("C15").Value = "dog/dog088.png"
$ws.Range("D15").Value = "langen"
$ws.Range("E15").Value = "dog"
$ws.Range("B16").Value = 45
$ws.Range("C16").Value = "flower/flower068.png"
$ws.Range("D16").Value = "strahlen"
$ws.Range("E16").Value = "flower"
$ws.Range("B17").Value = 43
$ws.Range("C17").Value = "flower/flower076.png"
$ws.Range("D17").Value = "posten"
$ws.Range("E17").Value = "flower"
$ws.Range("B18").Value = 87
$ws.Range("C18").Value = "flower/flower085.png"
$ws.Range("D18").Value = "laufen"
$ws.Range("E18").Value = "flower"
$ws.Range("B19").Value = 94
$ws.Range("C19").Value = "dog/dog087.png"
$ws.Range("D19").Value = "rücken"
$ws.Range("E19").Value = "dog"
$ws.Range("B20").Value = 23
$ws.Range("C20").Value = "flower/flower079.png"
$ws.Range("D20").Value = "schicken"
$ws.Range("E20").Value = "flower"
$ws.Range("B21").Value = 97
$ws.Range("C21").Value = "flower/flower073.png"
$ws.Range("D21").Value = "klappen"
$ws.Range("E21").Value = "flower"
$ws.Range("B22").Value = 81
$ws.Range("C22").Value = "flower/flower067.png"
$ws.Range("D22").Value = "krachen"
$ws.Range("E22").Value = "flower"
$ws.Range("B23").Value = 70
$ws.Range("C23").Value = "flower/flower066.png"
$ws.Range("D23").Value = "bitten"
$ws.Range("E23").Value = "flower"
$ws.Range("B24").Value = 92
$ws.Range("C24").Value = "flower/flower080.png"
$ws.Range("D24").Value = "formen"
$ws.Range("E24").Value = "flower"
$ws.Range("B25").Value = 106
$ws.Range("C25").Value = "dog/dog116.png"
$ws.Range("D25").Value = "hoffen"
$ws.Range("E25").Value = "dog"
$ws.Range("B26").Value = 105
$ws.Range("C26").Value = "dog/dog118.png"
$ws.Range("D26").Value = "spielen"
$ws.Range("E26").Value = "dog"
$ws.Range("B27").Value = 4
$ws.Range("C27").Value = "flower/flower094.png"
$ws.Range("D27").Value = "jubeln"
$ws.Range("E27").Value = "flower"
$ws.Range("B28").Value = 12
$ws.Range("C28").Value = "flower/flower070.png"
$ws.Range("D28").Value = "tauschen"
$ws.Range("E28").Value = "flower"
$ws.Range("B29").Value = 114
$ws.Range("C29").Value = "dog/dog110.png"
$ws.Range("D29").Value = "füllen"
$ws.Range("E29").Value = "dog"
$ws.Range("B30").Value = 66
$ws.Range("C30").Value = "dog/dog067.png"
$ws.Range("D30").Value = "runden"
$ws.Range("E30").Value = "dog"
$ws.Range("B31").Value = 95
$ws.Range("C31").Value = "dog/dog064.png"
$ws.Range("D31").Value = "scheitern"
$ws.Range("E31").Value = "dog"
$ws.Range("B32").Value = 2
$ws.Range("C32").Value = "dog/dog085.png"
$ws.Range("D32").Value = "hauen"
$ws.Range("E32").Value = "dog"
$ws.Range("B33").Value = 11
$ws.Range("C33").Value = "dog/dog095.png"
$ws.Range("D33").Value = "fesseln"
$ws.Range("E33").Value = "dog"
